# Applies the Refassist report regeneration (journal article -> conference paper)
# by locating each report line via Find and rewriting its Range.Text directly
# (Range.Text avoids Word's smart-quote autocorrect that Find.Execute's
# ReplaceWith parameter would otherwise trigger).
$d = $word.ActiveDocument
$vt = [char]11   # manual line break marker; becomes <w:br/> when written back out

# Change 1
$rng = $d.Content
$rng.Find.Execute("- Type detected: journal article") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- Type detected: conference paper" } else { Write-Output "WARNING: change 1 search text not found" }

# Change 2
$rng = $d.Content
$rng.Find.Execute("- DOI: 10.1016/j.jsv.2004.07.007") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- DOI: 10.1109/tpami.2018.2844175" } else { Write-Output "WARNING: change 2 search text not found" }

# Change 3
$rng = $d.Content
$rng.Find.Execute("- journal_abbrev: J. Sound Vib. → Missing journal name  (source: Unknown)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- journal_abbrev: MISSING → Missing journal name  (source: Unknown)" } else { Write-Output "WARNING: change 3 search text not found" }

# Change 4
$rng = $d.Content
$rng.Find.Execute("- authors: Holopainen T., Tenhunen A., Arkio A. → Timo Holopainen, A. Tenhunen, Antero Arkkio  (source: Unknown)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- authors: K. He, G. Gkioxari, P. Dollár, R. Girshick → Kaiming He, Georgia Gkioxari, Piotr Dollár, Ross Girshick  (source: Unknown)" } else { Write-Output "WARNING: change 4 search text not found" }

# Change 5
$rng = $d.Content
$rng.Find.Execute("- journal_name: MISSING → Journal of Sound and Vibration  (source: Unknown)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- journal_name: MISSING → IEEE Transactions on Pattern Analysis and Machine Intelligence  (source: Unknown)" } else { Write-Output "WARNING: change 5 search text not found" }

# Change 6
$rng = $d.Content
$rng.Find.Execute("- journal_abbrev: J. Sound Vib. → Journal of Sound and Vibration  (source: Unknown)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- journal_abbrev: MISSING → IEEE Trans. Pattern Anal. Mach. Intell.  (source: Unknown)" } else { Write-Output "WARNING: change 6 search text not found" }

# Change 7
$rng = $d.Content
$rng.Find.Execute("- issue: 3–5 → 3-5  (source: Unknown)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- volume: MISSING → 42  (source: Unknown)" } else { Write-Output "WARNING: change 7 search text not found" }

# Change 8
$rng = $d.Content
$rng.Find.Execute("- pages: 733–755 → 733-755  (source: Unknown)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- issue: MISSING → 2  (source: Unknown)" } else { Write-Output "WARNING: change 8 search text not found" }

# Change 9
$rng = $d.Content
$rng.Find.Execute("- doi: MISSING → 10.1016/j.jsv.2004.07.007  (source: Unknown)" + $vt + "- month: MISSING → 6  (source: Unknown)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- pages: 2961-2969 → 386-397  (source: Unknown)" + $vt + "- doi: MISSING → 10.1109/tpami.2018.2844175  (source: Unknown)" + $vt + "- year: 2017 → 2020  (source: Unknown)" + $vt + "- month: MISSING → 2  (source: Unknown)" } else { Write-Output "WARNING: change 9 search text not found" }

# Change 10
$rng = $d.Content
$rng.Find.Execute("- title: Electromechanical interaction in rotordynamics of cage induction motors  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- title: Mask R-CNN  (source: Not available)" } else { Write-Output "WARNING: change 10 search text not found" }

# Change 11
$rng = $d.Content
$rng.Find.Execute("- authors: Timo Holopainen, A. Tenhunen, Antero Arkkio  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- authors: Kaiming He, Georgia Gkioxari, Piotr Dollár, Ross Girshick  (source: Not available)" } else { Write-Output "WARNING: change 11 search text not found" }

# Change 12
$rng = $d.Content
$rng.Find.Execute("- journal_name: Journal of Sound and Vibration  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- journal_name: IEEE Transactions on Pattern Analysis and Machine Intelligence  (source: Not available)" } else { Write-Output "WARNING: change 12 search text not found" }

# Change 13
$rng = $d.Content
$rng.Find.Execute("- journal_abbrev: Journal of Sound and Vibration  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- journal_abbrev: IEEE Trans. Pattern Anal. Mach. Intell.  (source: Not available)" + $vt + "- conference_name: Proceedings of the IEEE International Conference on Computer Vision (ICCV)  (source: Not available)" } else { Write-Output "WARNING: change 13 search text not found" }

# Change 14
$rng = $d.Content
$rng.Find.Execute("- volume: 284  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- volume: 42  (source: Not available)" } else { Write-Output "WARNING: change 14 search text not found" }

# Change 15
$rng = $d.Content
$rng.Find.Execute("- issue: 3-5  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- issue: 2  (source: Not available)" } else { Write-Output "WARNING: change 15 search text not found" }

# Change 16
$rng = $d.Content
$rng.Find.Execute("- pages: 733-755  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- pages: 386-397  (source: Not available)" } else { Write-Output "WARNING: change 16 search text not found" }

# Change 17
$rng = $d.Content
$rng.Find.Execute("- year: 2005  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- year: 2020  (source: Not available)" } else { Write-Output "WARNING: change 17 search text not found" }

# Change 18
$rng = $d.Content
$rng.Find.Execute("- month: 6  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- month: 2  (source: Not available)" } else { Write-Output "WARNING: change 18 search text not found" }

# Change 19
$rng = $d.Content
$rng.Find.Execute("- doi: 10.1016/j.jsv.2004.07.007  (source: Not available)") | Out-Null
if ($rng.Find.Found) { $rng.Text = "- doi: 10.1109/tpami.2018.2844175  (source: Not available)" } else { Write-Output "WARNING: change 19 search text not found" }

# Change 20
$rng = $d.Content
$rng.Find.Execute("T. Holopainen, A. Tenhunen, and A. Arkkio, `"Electromechanical interaction in rotordynamics of cage induction motors,`" *Journal of Sound and Vibration*, vol. 284, no. 3-5, pp. 733–755, Jun. 2005, https://doi.org/10.1016/j.jsv.2004.07.007.") | Out-Null
if ($rng.Find.Found) { $rng.Text = "K. He, G. Gkioxari, P. Dollár, and R. Girshick, `"Mask R-CNN,`" in *Proceedings of the IEEE International Conference on Computer Vision (ICCV)*, pp. 386-397, Feb. 2020, https://doi.org/10.1109/tpami.2018.2844175." } else { Write-Output "WARNING: change 20 search text not found" }

